$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Femacal de La Calera" / Espinaca.
# In the saved file this shows up as a new row inserted at row 193, which
# pushes the existing rows 193-226 down to 194-227 (dimension grows from
# A1:R226 to A1:R227). Replicate that with a real row insert so every
# shifted row (and its formatting) moves down exactly one position.
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row with the new week's data.
$ws.Range("A193").Value = 3
$ws.Range("B193").Value = "Femacal de La Calera"
$ws.Range("C193").Value = "Coquimbo"
$ws.Range("D193").Value = 44522
$ws.Range("E193").Value = 5
$ws.Range("F193").Value = 100112012
$ws.Range("G193").Value = "Espinaca"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 270
$ws.Range("K193").Value = 2500
$ws.Range("L193").Value = 2800
$ws.Range("M193").Value = 2667
$ws.Range("N193").Value = '$/docena de atados (3 kilos)'
$ws.Range("O193").Value = "Provincia de Quillota"
$ws.Range("P193").Value = 889
$ws.Range("Q193").Value = 3
$ws.Range("R193").Value = "Hortaliza"
